$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col4a1"
$ws.Range("C2").Value = "Itgb8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 161.279784
$ws.Range("H2").Value = 483.839352
$ws.Range("I2").Value = 0.3023989599621841
$ws.Range("J2").Value = 0.3023989599621841
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1126243333333333
$ws.Range("N2").Value = 0.337873
$ws.Range("O2").Value = 0.01082936903163217
$ws.Range("P2").Value = 0.01082936903163217
$ws.Range("Q2").Value = 18.164028153144
$ws.Range("R2").Value = 163.476253378296
$ws.Range("S2").Value = 0.003274789932212252
$ws.Range("T2").Value = 0.003274789932212252

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col4a1"
$ws.Range("C3").Value = "Itgb8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 161.279784
$ws.Range("H3").Value = 483.839352
$ws.Range("I3").Value = 0.3023989599621841
$ws.Range("J3").Value = 0.3023989599621841
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.323421
$ws.Range("N3").Value = 9.970263
$ws.Range("O3").Value = 0.3195628457125252
$ws.Range("P3").Value = 0.3195628457125252
$ws.Range("Q3").Value = 536.0006210210639
$ws.Range("R3").Value = 4824.005589189575
$ws.Range("S3").Value = 0.09663547218602353
$ws.Range("T3").Value = 0.09663547218602353

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col4a1"
$ws.Range("C4").Value = "Itgb8"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 161.279784
$ws.Range("H4").Value = 483.839352
$ws.Range("I4").Value = 0.3023989599621841
$ws.Range("J4").Value = 0.3023989599621841
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.963852666666667
$ws.Range("N4").Value = 20.891558
$ws.Range("O4").Value = 0.6696077852558425
$ws.Range("P4").Value = 0.6696077852558425
$ws.Range("Q4").Value = 1123.128653887824
$ws.Range("R4").Value = 10108.15788499041
$ws.Range("S4").Value = 0.2024886978439483
$ws.Range("T4").Value = 0.2024886978439483

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col4a1"
$ws.Range("C5").Value = "Itgb8"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 288.7700093333333
$ws.Range("H5").Value = 866.3100279999999
$ws.Range("I5").Value = 0.541442630470476
$ws.Range("J5").Value = 0.5414426304704759
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1126243333333333
$ws.Range("N5").Value = 0.337873
$ws.Range("O5").Value = 0.01082936903163217
$ws.Range("P5").Value = 0.01082936903163217
$ws.Range("Q5").Value = 32.5225297878271
$ws.Range("R5").Value = 292.702768090444
$ws.Range("S5").Value = 0.005863482054822433
$ws.Range("T5").Value = 0.005863482054822431

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col4a1"
$ws.Range("C6").Value = "Itgb8"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 288.7700093333333
$ws.Range("H6").Value = 866.3100279999999
$ws.Range("I6").Value = 0.541442630470476
$ws.Range("J6").Value = 0.5414426304704759
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.323421
$ws.Range("N6").Value = 9.970263
$ws.Range("O6").Value = 0.3195628457125252
$ws.Range("P6").Value = 0.3195628457125252
$ws.Range("Q6").Value = 959.7043131885957
$ws.Range("R6").Value = 8637.338818697363
$ws.Range("S6").Value = 0.1730249477832206
$ws.Range("T6").Value = 0.1730249477832205

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col4a1"
$ws.Range("C7").Value = "Itgb8"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 288.7700093333333
$ws.Range("H7").Value = 866.3100279999999
$ws.Range("I7").Value = 0.541442630470476
$ws.Range("J7").Value = 0.5414426304704759
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.963852666666667
$ws.Range("N7").Value = 20.891558
$ws.Range("O7").Value = 0.6696077852558425
$ws.Range("P7").Value = 0.6696077852558425
$ws.Range("Q7").Value = 2010.951799549291
$ws.Range("R7").Value = 18098.56619594362
$ws.Range("S7").Value = 0.362554200632433
$ws.Range("T7").Value = 0.3625542006324329

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col4a1"
$ws.Range("C8").Value = "Itgb8"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 83.28466000000002
$ws.Range("H8").Value = 249.85398
$ws.Range("I8").Value = 0.15615840956734
$ws.Range("J8").Value = 0.15615840956734
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1126243333333333
$ws.Range("N8").Value = 0.337873
$ws.Range("O8").Value = 0.01082936903163217
$ws.Range("P8").Value = 0.01082936903163217
$ws.Range("Q8").Value = 9.379879309393335
$ws.Range("R8").Value = 84.41891378454001
$ws.Range("S8").Value = 0.001691097044597484
$ws.Range("T8").Value = 0.001691097044597484

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col4a1"
$ws.Range("C9").Value = "Itgb8"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 83.28466000000002
$ws.Range("H9").Value = 249.85398
$ws.Range("I9").Value = 0.15615840956734
$ws.Range("J9").Value = 0.15615840956734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.323421
$ws.Range("N9").Value = 9.970263
$ws.Range("O9").Value = 0.3195628457125252
$ws.Range("P9").Value = 0.3195628457125252
$ws.Range("Q9").Value = 276.78998802186
$ws.Range("R9").Value = 2491.10989219674
$ws.Range("S9").Value = 0.04990242574328118
$ws.Range("T9").Value = 0.04990242574328118

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col4a1"
$ws.Range("C10").Value = "Itgb8"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 83.28466000000002
$ws.Range("H10").Value = 249.85398
$ws.Range("I10").Value = 0.15615840956734
$ws.Range("J10").Value = 0.15615840956734
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.963852666666667
$ws.Range("N10").Value = 20.891558
$ws.Range("O10").Value = 0.6696077852558425
$ws.Range("P10").Value = 0.6696077852558425
$ws.Range("Q10").Value = 579.9821016334269
$ws.Range("R10").Value = 5219.838914700841
$ws.Range("S10").Value = 0.1045648867794613
$ws.Range("T10").Value = 0.1045648867794613
